# Translate the Sprache (language) column from German into English, and
# add a new "Language Familiy" column (D) that classifies each language
# into its language family (Romance / Germanic).
#
# The workbook currently has:
#   A: Club, B: Land, C: Sprache  (values: Deutsch / Italienisch / Französisch)
# After this script it should have:
#   A: Club, B: Land, C: Sprache (values: German / Italian / French)
#   D: Language Familiy          (values: Germanic / Romance)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 162

# --- 1) Translate column C (Sprache) values to English -------------------
# Do this one language at a time so that every reference to the old
# German-language shared string is removed before the next one is
# introduced (keeps the shared-string table compact / in a sane order).
for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Value2 -eq "Deutsch") {
        $ws.Cells.Item($r, 3).Value2 = "German"
    }
}
for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Value2 -eq "Italienisch") {
        $ws.Cells.Item($r, 3).Value2 = "Italian"
    }
}
for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Value2 -eq "Französisch") {
        $ws.Cells.Item($r, 3).Value2 = "French"
    }
}

# --- 2) Add the new "Language Familiy" column (D) -------------------------
$ws.Cells.Item(1, 4).Value2 = "Language Familiy"

for ($r = 2; $r -le $lastRow; $r++) {
    $lang = $ws.Cells.Item($r, 3).Value2
    if ($lang -eq "German") {
        $ws.Cells.Item($r, 4).Value2 = "Germanic"
    }
    else {
        # Italian and French are both Romance languages
        $ws.Cells.Item($r, 4).Value2 = "Romance"
    }
}

# --- 3) Column width for the new column -----------------------------------
$ws.Columns.Item(4).ColumnWidth = 14.7

# --- 4) Update the selection (mirrors the author re-selecting a cell) -----
$ws.Range("J23").Select() | Out-Null

# --- 5) Refresh the AutoFilter to cover the new column ---------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:D162").AutoFilter() | Out-Null

# --- 6) Update the _FilterDatabase defined name to match the new range ----
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=swiss_teams_lamguage!`$A`$1:`$D`$162"
    }
}
